{"js": "// Update the date line and the 25 three-digit-by-one-digit multiplication\n// answers in the practice-sheet table to the new values from the next day's\n// worksheet (2024-08-18 Sunday).\nconst replacements = [\n  [\"2024-08-17 Saturday\", \"2024-08-18 Sunday\"],\n  [\"677\u00d77=4739\", \"521\u00d73=1563\"],\n  [\"864\u00d74=3456\", \"730\u00d73=2190\"],\n  [\"564\u00d79=5076\", \"812\u00d76=4872\"],\n  [\"147\u00d74=588\", \"143\u00d75=715\"],\n  [\"757\u00d78=6056\", \"954\u00d75=4770\"],\n  [\"521\u00d74=2084\", \"191\u00d72=382\"],\n  [\"869\u00d73=2607\", \"736\u00d72=1472\"],\n  [\"552\u00d75=2760\", \"705\u00d77=4935\"],\n  [\"532\u00d75=2660\", \"337\u00d74=1348\"],\n  [\"840\u00d78=6720\", \"821\u00d72=1642\"],\n  [\"998\u00d72=1996\", \"451\u00d78=3608\"],\n  [\"662\u00d77=4634\", \"228\u00d79=2052\"],\n  [\"258\u00d72=516\", \"829\u00d78=6632\"],\n  [\"300\u00d79=2700\", \"229\u00d73=687\"],\n  [\"444\u00d79=3996\", \"798\u00d78=6384\"],\n  [\"775\u00d72=1550\", \"497\u00d74=1988\"],\n  [\"840\u00d72=1680\", \"317\u00d78=2536\"],\n  [\"187\u00d79=1683\", \"441\u00d73=1323\"],\n  [\"492\u00d77=3444\", \"725\u00d73=2175\"],\n  [\"719\u00d77=5033\", \"382\u00d78=3056\"],\n  [\"991\u00d75=4955\", \"284\u00d77=1988\"],\n  [\"292\u00d77=2044\", \"233\u00d77=1631\"],\n  [\"969\u00d73=2907\", \"959\u00d78=7672\"],\n  [\"149\u00d74=596\", \"443\u00d79=3987\"],\n  [\"677\u00d78=5416\", \"869\u00d74=3476\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 three-digit-by-one-digit multiplication\n# answers in the practice-sheet table to the new values from the next day's\n# worksheet (2024-08-18 Sunday).\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-08-17 Saturday\", \"2024-08-18 Sunday\"),\n    @(\"677\u00d77=4739\", \"521\u00d73=1563\"),\n    @(\"864\u00d74=3456\", \"730\u00d73=2190\"),\n    @(\"564\u00d79=5076\", \"812\u00d76=4872\"),\n    @(\"147\u00d74=588\", \"143\u00d75=715\"),\n    @(\"757\u00d78=6056\", \"954\u00d75=4770\"),\n    @(\"521\u00d74=2084\", \"191\u00d72=382\"),\n    @(\"869\u00d73=2607\", \"736\u00d72=1472\"),\n    @(\"552\u00d75=2760\", \"705\u00d77=4935\"),\n    @(\"532\u00d75=2660\", \"337\u00d74=1348\"),\n    @(\"840\u00d78=6720\", \"821\u00d72=1642\"),\n    @(\"998\u00d72=1996\", \"451\u00d78=3608\"),\n    @(\"662\u00d77=4634\", \"228\u00d79=2052\"),\n    @(\"258\u00d72=516\", \"829\u00d78=6632\"),\n    @(\"300\u00d79=2700\", \"229\u00d73=687\"),\n    @(\"444\u00d79=3996\", \"798\u00d78=6384\"),\n    @(\"775\u00d72=1550\", \"497\u00d74=1988\"),\n    @(\"840\u00d72=1680\", \"317\u00d78=2536\"),\n    @(\"187\u00d79=1683\", \"441\u00d73=1323\"),\n    @(\"492\u00d77=3444\", \"725\u00d73=2175\"),\n    @(\"719\u00d77=5033\", \"382\u00d78=3056\"),\n    @(\"991\u00d75=4955\", \"284\u00d77=1988\"),\n    @(\"292\u00d77=2044\", \"233\u00d77=1631\"),\n    @(\"969\u00d73=2907\", \"959\u00d78=7672\"),\n    @(\"149\u00d74=596\", \"443\u00d79=3987\"),\n    @(\"677\u00d78=5416\", \"869\u00d74=3476\")\n)\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n}\n"}
